$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 53
$ws.Range("D6").Value = 90.59999999999999

$ws.Range("D7").Value = 94.3

$ws.Range("A8").Value = "MediaTek MT7921 Wi-Fi 6 802.11ax PCIe Adapter - 3.0.1.1327"
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 418
$ws.Range("D8").Value = 95.09999999999999

$ws.Range("A9").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.100.1.1"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 96.40000000000001

$ws.Range("B10").Value = 2

$ws.Range("A11").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.200.2.1"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 97.5

$ws.Range("A12").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.70.2.3"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 97.7

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.80.0.7"
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 97.8

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.70.0.6"
$ws.Range("C14").Value = 24
$ws.Range("D14").Value = 98.09999999999999

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.120.0.3"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 98.3

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.170.0.3"
$ws.Range("B16").Value = 12
$ws.Range("C16").Value = 197
$ws.Range("D16").Value = 98.3

$ws.Range("A17").Value = "MediaTek MT7921 Wi-Fi 6 802.11ax PCIe Adapter - 3.0.1.1303"
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 35

$ws.Range("D19").Value = 98.7

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.0.3"
$ws.Range("B21").Value = 3
$ws.Range("C21").Value = 115

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.20.1.1"
$ws.Range("B22").Value = 11
$ws.Range("C22").Value = 54
$ws.Range("D22").Value = 98.90000000000001

$ws.Range("B23").Value = 64
$ws.Range("C23").Value = 1637
